# Apply the deck-level theme re-colour (Integral -> Office) to the
# presentation's (slide master) theme, and fix up the table style on the
# table that still references the old, document-local "Table_0" style so
# that it uses the standard built-in table style instead.

$p = $ppt.ActivePresentation

# --- 1) Re-colour the presentation theme (ppt/theme/theme1.xml) --------
# Target colour scheme ("Office"), expressed as COM RGB() values
# (0xBBGGRR packed long, i.e. R + G*256 + B*65536).
$m = $p.SlideMaster
$cs = $m.Theme.ThemeColorScheme

$cs.Item(1).RGB  = 0            # dk1      000000
$cs.Item(2).RGB  = 16777215     # lt1      FFFFFF
$cs.Item(3).RGB  = 6968388      # dk2      44546A
$cs.Item(4).RGB  = 15132391     # lt2      E7E6E6
$cs.Item(5).RGB  = 13998939     # accent1  5B9BD5
$cs.Item(6).RGB  = 3243501      # accent2  ED7D31
$cs.Item(7).RGB  = 10855845     # accent3  A5A5A5
$cs.Item(8).RGB  = 49407        # accent4  FFC000
$cs.Item(9).RGB  = 12874308     # accent5  4472C4
$cs.Item(10).RGB = 4697456      # accent6  70AD47
$cs.Item(11).RGB = 12673797     # hlink    0563C1
$cs.Item(12).RGB = 7491477      # folHlink 954F72

# --- 2) Fix the table style on slide 16 ---------------------------------
$s = $p.Slides.Item(16)
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.HasTable) {
        $tbl = $shp.Table
        $tbl.ApplyStyle("{41020D70-EABD-4F68-AE93-A1FAF7029921}")
    }
}
